$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.67
$ws.Range("I3").Value = 6
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("W3").Value = 5.5
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7
$ws.Range("AD3").Value = 6.5
$ws.Range("AL3").Value = 51
$ws.Range("AP3").Value = 23
$ws.Range("AQ3").Value = 29
$ws.Range("AR3").Value = 51
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 9.5
$ws.Range("AZ3").Value = 126
$ws.Range("BA3").Value = 151
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.9
$ws.Range("R5").Value = 1.95
$ws.Range("G6").Value = 2.45
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3.25
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 1.53
$ws.Range("Y6").Value = 10
$ws.Range("AJ6").Value = 34
$ws.Range("AM6").Value = 900
$ws.Range("AZ6").Value = 67
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 3.4
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1.85
$ws.Range("G9").Value = 2.35
$ws.Range("I9").Value = 2.75
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 3.25
$ws.Range("AI9").Value = 10
$ws.Range("AU9").Value = 7
$ws.Range("G14").Value = 1.91
$ws.Range("I14").Value = 3.75
$ws.Range("J14").Value = 2.75
$ws.Range("L14").Value = 4.75
$ws.Range("Q14").Value = 2.35
$ws.Range("R14").Value = 1.57
$ws.Range("S14").Value = 1.53
$ws.Range("T14").Value = 2.38
$ws.Range("U14").Value = 2.1
$ws.Range("V14").Value = 1.67
$ws.Range("X14").Value = 8
$ws.Range("Z14").Value = 17
$ws.Range("AG14").Value = 9
$ws.Range("AH14").Value = 19
$ws.Range("AI14").Value = 15
$ws.Range("AK14").Value = 41
$ws.Range("AN14").Value = 3.75
$ws.Range("AO14").Value = 11
$ws.Range("AT14").Value = 2.38
$ws.Range("AX14").Value = 23
$ws.Range("AZ14").Value = 81
$ws.Range("BA14").Value = 126
$ws.Range("BB14").Value = 351
$ws.Range("H16").Value = 3.2
$ws.Range("I16").Value = 2.75
$ws.Range("J16").Value = 3.1
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 3.5
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 1.33
$ws.Range("P16").Value = 3.25
$ws.Range("Q16").Value = 2.08
$ws.Range("R16").Value = 1.73
$ws.Range("U16").Value = 1.8
$ws.Range("V16").Value = 1.91
$ws.Range("W16").Value = 8
$ws.Range("X16").Value = 12
$ws.Range("AB16").Value = 29
$ws.Range("AC16").Value = 9
$ws.Range("AD16").Value = 6.5
$ws.Range("AG16").Value = 8.5
$ws.Range("AK16").Value = 23
$ws.Range("AM16").Value = 251
$ws.Range("AO16").Value = 13
$ws.Range("AP16").Value = 23
$ws.Range("AQ16").Value = 41
$ws.Range("AS16").Value = 151
$ws.Range("AU16").Value = 8
$ws.Range("AY16").Value = 26
$ws.Range("Q20").Value = 2.3
$ws.Range("R20").Value = 1.6
$ws.Range("G22").Value = 3.4
$ws.Range("I22").Value = 2.2
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 3
$ws.Range("U22").Value = 1.91
$ws.Range("V22").Value = 1.8
$ws.Range("W22").Value = 9
$ws.Range("Y22").Value = 13
$ws.Range("Z22").Value = 41
$ws.Range("AC22").Value = 7.5
$ws.Range("AE22").Value = 17
$ws.Range("AG22").Value = 6.5
$ws.Range("AH22").Value = 9.5
$ws.Range("AO22").Value = 21
$ws.Range("AQ22").Value = 67
$ws.Range("AR22").Value = 101
$ws.Range("AV22").Value = 67
$ws.Range("AW22").Value = 4
$ws.Range("G23").Value = 2.6
$ws.Range("I23").Value = 2.7
$ws.Range("J23").Value = 3.4
$ws.Range("L23").Value = 3.5
$ws.Range("W23").Value = 7.5
$ws.Range("X23").Value = 12
$ws.Range("Z23").Value = 26
$ws.Range("AH23").Value = 12
$ws.Range("AJ23").Value = 26
$ws.Range("AK23").Value = 23
$ws.Range("AL23").Value = 34
$ws.Range("AV23").Value = 51
$ws.Range("BB23").Value = 201
$ws.Range("G24").Value = 2.05
$ws.Range("I24").Value = 4
$ws.Range("L24").Value = 4.33
$ws.Range("M24").Value = 1.06
$ws.Range("N24").Value = 10
$ws.Range("O24").Value = 1.33
$ws.Range("P24").Value = 3.25
$ws.Range("Q24").Value = 2.05
$ws.Range("R24").Value = 1.75
$ws.Range("X24").Value = 9
$ws.Range("Z24").Value = 17
$ws.Range("AA24").Value = 17
$ws.Range("AC24").Value = 8
$ws.Range("AK24").Value = 34
$ws.Range("AO24").Value = 11
$ws.Range("G35").Value = 1.87
$ws.Range("H35").Value = 3.3
$ws.Range("I35").Value = 4
$ws.Range("J35").Value = 2.42
$ws.Range("L35").Value = 4.3
$ws.Range("O35").Value = 1.31
$ws.Range("P35").Value = 2.87
$ws.Range("Q35").Value = 1.91
$ws.Range("S35").Value = 1.4
$ws.Range("T35").Value = 2.52
$ws.Range("U35").Value = 1.78
$ws.Range("V35").Value = 1.83
$ws.Range("W35").Value = 6.6
$ws.Range("X35").Value = 8.5
$ws.Range("Z35").Value = 16
$ws.Range("AA35").Value = 15.5
$ws.Range("AB35").Value = 28
$ws.Range("AD35").Value = 6.4
$ws.Range("AE35").Value = 15
$ws.Range("AF35").Value = 75
$ws.Range("AG35").Value = 11
$ws.Range("AI35").Value = 13
$ws.Range("AK35").Value = 37
$ws.Range("AL35").Value = 45
$ws.Range("AM35").Value = 600
$ws.Range("AN35").Value = 3.7
$ws.Range("AO35").Value = 9.25
$ws.Range("AP35").Value = 17.5
$ws.Range("AS35").Value = 200
$ws.Range("AT35").Value = 2.57
$ws.Range("AU35").Value = 7
$ws.Range("AV35").Value = 60
$ws.Range("AX35").Value = 22
$ws.Range("AY35").Value = 27
$ws.Range("BB35").Value = 350
$ws.Range("G38").Value = 10.75
$ws.Range("H38").Value = 5
$ws.Range("J38").Value = 9
$ws.Range("K38").Value = 2.45
$ws.Range("L38").Value = 1.7
$ws.Range("P38").Value = 3.8
$ws.Range("U38").Value = 2.32
$ws.Range("V38").Value = 1.55
$ws.Range("W38").Value = 26
$ws.Range("X38").Value = 90
$ws.Range("Y38").Value = 37
$ws.Range("Z38").Value = 450
$ws.Range("AD38").Value = 10.75
$ws.Range("AI38").Value = 9.25
$ws.Range("AJ38").Value = 6.9
$ws.Range("AK38").Value = 11.5
$ws.Range("AN38").Value = 11
$ws.Range("AO38").Value = 75
$ws.Range("AP38").Value = 70
$ws.Range("H39").Value = 4.45
$ws.Range("I39").Value = 6.1
$ws.Range("J39").Value = 1.88
$ws.Range("L39").Value = 5.7
$ws.Range("P39").Value = 4.5
$ws.Range("S39").Value = 1.28
$ws.Range("T39").Value = 3.35
$ws.Range("X39").Value = 7.8
$ws.Range("Z39").Value = 10
$ws.Range("AG39").Value = 22
$ws.Range("AI39").Value = 19.5
$ws.Range("AP39").Value = 14.5
$ws.Range("AQ39").Value = 17.5
$ws.Range("AR39").Value = 40
$ws.Range("AT39").Value = 3.35
$ws.Range("AW39").Value = 7.9
$ws.Range("AX39").Value = 32
$ws.Range("AY39").Value = 32
$ws.Range("BB39").Value = 350
